$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: JEMAA HORMI / B219321, rappel only (MT brut/taxe rappel = 6000) ---
$ws.Range("A2").Value = "001/RRR"
$ws.Range("B2").Value = "Direction régionale"
$ws.Range("C2").Value = "B219321"
$ws.Range("D2").Value = "JEMAA HORMI"
$ws.Range("E2").Value = "non"
$ws.Range("F2").Value = "mensuelle"
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = "--"
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = "--"
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 6000
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = "--"
$ws.Range("O2").Value = 6000

# --- Row 3: NOUBAIL MOHAMMED / IR801997, rappel only (MT brut/taxe rappel = 6000) ---
$ws.Range("A3").Value = "001/RRR"
$ws.Range("B3").Value = "Direction régionale"
$ws.Range("C3").Value = "IR801997"
$ws.Range("D3").Value = "NOUBAIL MOHAMMED"
$ws.Range("E3").Value = "non"
$ws.Range("F3").Value = "mensuelle"
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = "--"
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = "--"
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 6000
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = "--"
$ws.Range("O3").Value = 6000

# --- Row 4: JEMAA HORMI / B219321, regular monthly rent line (MT brut loyer = 1000) ---
$ws.Range("A4").Value = "001/RRR"
$ws.Range("B4").Value = "Direction régionale"
$ws.Range("C4").Value = "B219321"
$ws.Range("D4").Value = "JEMAA HORMI"
$ws.Range("E4").Value = "non"
$ws.Range("F4").Value = "mensuelle"
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 1000
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 1000

# --- Row 5: NOUBAIL MOHAMMED / IR801997, regular monthly rent line (MT brut loyer = 1000) ---
$ws.Range("A5").Value = "001/RRR"
$ws.Range("B5").Value = "Direction régionale"
$ws.Range("C5").Value = "IR801997"
$ws.Range("D5").Value = "NOUBAIL MOHAMMED"
$ws.Range("E5").Value = "non"
$ws.Range("F5").Value = "mensuelle"
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 1000
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = 1000

# --- Row 6: totals row (was row 3) ---
$ws.Range("A6").Value = " "
$ws.Range("B6").Value = " "
$ws.Range("C6").Value = " "
$ws.Range("D6").Value = " "
$ws.Range("E6").Value = " "
$ws.Range("F6").Value = " "
$ws.Range("G6").Value = " "
$ws.Range("H6").Value = 2000
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 12000
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 0
$ws.Range("O6").Value = 14000
